# PPDM Website Cookies Data - add "Vendors" column with First/Third party info,
# and add a hyperlink for the Bright Bridge website URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before the current column E (Performance Cookies),
# this shifts the existing E..H columns to F..I and copies formatting/width.
$ws.Columns("E:E").Insert()

# New column header + the "One Trust" vendor detail cell.
$ws.Range("E3").Value = "Vendors"
$ws.Range("E4").Value = "First Party - Optanon Consent  & AW SalB.. Third Party-  _GRECAPTCHA"

# Give the new column the same width as column D.
$ws.Range("E1").ColumnWidth = 54.43

# Add hyperlink for the "Bright Bridge" website url cell.
$ws.Hyperlinks.Add($ws.Range("C5"), "https://brightbridgesolutions.com/")

# Update the active selection to C5 (matches the saved view state).
$ws.Range("C5").Select()
